$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.09665433333333333
$ws.Range("H2").Value = 0.289963
$ws.Range("I2").Value = 0.0006230336790718351
$ws.Range("J2").Value = 0.0006230336790718351
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1007803333333333
$ws.Range("N2").Value = 0.302341
$ws.Range("O2").Value = 0.1035761647483865
$ws.Range("P2").Value = 0.1035761647483865
$ws.Range("Q2").Value = 0.009740855931444446
$ws.Range("R2").Value = 0.08766770338299999
$ws.Range("S2").Value = 0.00006453143898733777
$ws.Range("T2").Value = 0.00006453143898733778

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.09665433333333333
$ws.Range("H3").Value = 0.289963
$ws.Range("I3").Value = 0.0006230336790718351
$ws.Range("J3").Value = 0.0006230336790718351
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.623012
$ws.Range("N3").Value = 1.869036
$ws.Range("O3").Value = 0.6402954963325033
$ws.Range("P3").Value = 0.6402954963325033
$ws.Range("Q3").Value = 0.06021680951866666
$ws.Range("R3").Value = 0.5419512856679999
$ws.Range("S3").Value = 0.0003989256587731662
$ws.Range("T3").Value = 0.0003989256587731662

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.09665433333333333
$ws.Range("H4").Value = 0.289963
$ws.Range("I4").Value = 0.0006230336790718351
$ws.Range("J4").Value = 0.0006230336790718351
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2492146666666667
$ws.Range("N4").Value = 0.747644
$ws.Range("O4").Value = 0.2561283389191102
$ws.Range("P4").Value = 0.2561283389191102
$ws.Range("Q4").Value = 0.02408767746355555
$ws.Range("R4").Value = 0.216789097172
$ws.Range("S4").Value = 0.0001595765813113311
$ws.Range("T4").Value = 0.0001595765813113311

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 154.8642143333334
$ws.Range("H5").Value = 464.5926430000001
$ws.Range("I5").Value = 0.9982544794956518
$ws.Range("J5").Value = 0.9982544794956519
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1007803333333333
$ws.Range("N5").Value = 0.302341
$ws.Range("O5").Value = 0.1035761647483865
$ws.Range("P5").Value = 0.1035761647483865
$ws.Range("Q5").Value = 15.60726714191812
$ws.Range("R5").Value = 140.465404277263
$ws.Range("S5").Value = 0.1033953704290565
$ws.Range("T5").Value = 0.1033953704290565

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 154.8642143333334
$ws.Range("H6").Value = 464.5926430000001
$ws.Range("I6").Value = 0.9982544794956518
$ws.Range("J6").Value = 0.9982544794956519
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.623012
$ws.Range("N6").Value = 1.869036
$ws.Range("O6").Value = 0.6402954963325033
$ws.Range("P6").Value = 0.6402954963325033
$ws.Range("Q6").Value = 96.48226390023869
$ws.Range("R6").Value = 868.3403751021481
$ws.Range("S6").Value = 0.6391778474148131
$ws.Range("T6").Value = 0.6391778474148132

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 154.8642143333334
$ws.Range("H7").Value = 464.5926430000001
$ws.Range("I7").Value = 0.9982544794956518
$ws.Range("J7").Value = 0.9982544794956519
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2492146666666667
$ws.Range("N7").Value = 0.747644
$ws.Range("O7").Value = 0.2561283389191102
$ws.Range("P7").Value = 0.2561283389191102
$ws.Range("Q7").Value = 38.5944335536769
$ws.Range("R7").Value = 347.349901983092
$ws.Range("S7").Value = 0.2556812616517822
$ws.Range("T7").Value = 0.2556812616517823

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.174137
$ws.Range("H8").Value = 0.522411
$ws.Range("I8").Value = 0.001122486825276316
$ws.Range("J8").Value = 0.001122486825276316
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1007803333333333
$ws.Range("N8").Value = 0.302341
$ws.Range("O8").Value = 0.1035761647483865
$ws.Range("P8").Value = 0.1035761647483865
$ws.Range("Q8").Value = 0.01754958490566667
$ws.Range("R8").Value = 0.157946264151
$ws.Range("S8").Value = 0.0001162628803427131
$ws.Range("T8").Value = 0.0001162628803427131

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.174137
$ws.Range("H9").Value = 0.522411
$ws.Range("I9").Value = 0.001122486825276316
$ws.Range("J9").Value = 0.001122486825276316
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.623012
$ws.Range("N9").Value = 1.869036
$ws.Range("O9").Value = 0.6402954963325033
$ws.Range("P9").Value = 0.6402954963325033
$ws.Range("Q9").Value = 0.108489440644
$ws.Range("R9").Value = 0.9764049657959999
$ws.Range("S9").Value = 0.0007187232589169945
$ws.Range("T9").Value = 0.0007187232589169946

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.174137
$ws.Range("H10").Value = 0.522411
$ws.Range("I10").Value = 0.001122486825276316
$ws.Range("J10").Value = 0.001122486825276316
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.2492146666666667
$ws.Range("N10").Value = 0.747644
$ws.Range("O10").Value = 0.2561283389191102
$ws.Range("P10").Value = 0.2561283389191102
$ws.Range("Q10").Value = 0.04339749440933333
$ws.Range("R10").Value = 0.3905774496839999
$ws.Range("S10").Value = 0.0002875006860166083
$ws.Range("T10").Value = 0.0002875006860166083

